$d = $word.ActiveDocument

# 1. "Permet de rajouter une " -> "Permet d’ajouter une "
$d.Content.Find.Execute("Permet de rajouter une", $true, $false, $false, $false, $false, $true, 1, $false, "Permet d’ajouter une", 2) | Out-Null

# 2. "Administrateur désir ajouté un élément." -> "Administrateur désire ajouter un élément."
$d.Content.Find.Execute("Administrateur désir ajouté un élément.", $true, $false, $false, $false, $false, $true, 1, $false, "Administrateur désire ajouter un élément.", 2) | Out-Null

# 3. " qu’il désir." -> " qu’il désire."
$d.Content.Find.Execute("qu’il désir.", $true, $false, $false, $false, $false, $true, 1, $false, "qu’il désire.", 2) | Out-Null

# 4. "ppuis sur le bouton copier" -> "ppuie sur le bouton copier"
$d.Content.Find.Execute("ppuis sur le bouton copier", $true, $false, $false, $false, $false, $true, 1, $false, "ppuie sur le bouton copier", 2) | Out-Null

# 5. " et retourne le ID de l’élément" -> " et retourne l’ID de l’élément"
$d.Content.Find.Execute("et retourne le ID de l’élément", $true, $false, $false, $false, $false, $true, 1, $false, "et retourne l’ID de l’élément", 2) | Out-Null

# 6. "    Retourne l’ ID de l’élément" -> "    Retourne l’ID de l’élément" (also drops proofErr + moves bookmark away)
$d.Content.Find.Execute("Retourne l’ ID de l’élément", $true, $false, $false, $false, $false, $true, 1, $false, "Retourne l’ID de l’élément", 2) | Out-Null

# 7. "L’administrateur saisi le " -> "L’administrateur saisit le "
$d.Content.Find.Execute("L’administrateur saisi le", $true, $false, $false, $false, $false, $true, 1, $false, "L’administrateur saisit le", 2) | Out-Null

# 8. " saisie les informations qu’il veut modifier" -> " saisit les informations qu’il veut modifier"
$d.Content.Find.Execute("saisie les informations qu’il veut modifier", $true, $false, $false, $false, $false, $true, 1, $false, "saisit les informations qu’il veut modifier", 2) | Out-Null

# 9. "Le système montre u" -> "Le système affiche u"
$d.Content.Find.Execute("Le système montre u", $true, $false, $false, $false, $false, $true, 1, $false, "Le système affiche u", 2) | Out-Null

# 10. "» peut être déplacé" -> "» peut être déplacée"
$d.Content.Find.Execute("peut être déplacé", $true, $false, $false, $false, $false, $true, 1, $false, "peut être déplacée", 2) | Out-Null
